$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(8, 1).Value = "VET_DEP :"
$ws.Cells.Item(9, 1).Value = "И.О. заведущего Марьинским ветучастком VET_DEP_______________     ___________       VET_CEO"
$ws.Cells.Item(10, 1).Value = "  Специалист __________________________________________________  VET_DOC                     "

$ws.Range("A11").Select() | Out-Null

